# "update db y score"
# - Adds new cells on the "Score" sheet (H21, J21, J22, J23, K23) carrying
#   helper lookup labels.
# - Adds a brand-new "Hoja1" worksheet (placed after "Score", and left as the
#   active sheet) containing a small scoring-rules table (rows 2-11), a
#   second header/summary table (rows 17-18) with a SUM formula, and the
#   currency/percent/date formatting that goes with it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Score")

# ---------------------------------------------------------------------
# 1. Score sheet: a handful of new helper cells further out in the sheet
# ---------------------------------------------------------------------
$ws1.Range("H21").Value = "regla_edad"
$ws1.Range("J21").Value = "socres"
$ws1.Range("J22").Value = "entidad_id"
$ws1.Range("J23").Value = "score"
$ws1.Range("K23").Value = "num"

# Restore the prior selection on the Score sheet (it is no longer the
# active sheet once Hoja1 is added/activated below).
$ws1.Range("B44").Select()

# ---------------------------------------------------------------------
# 2. New worksheet "Hoja1", inserted right after "Score"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja1"

# --- Rules table header (row 2) ---
$ws2.Range("B2").Value = "CIA"
$ws2.Range("C2").Value = "REGLA"
$ws2.Range("D2").Value = "TIPO DE REGLA"
$ws2.Range("E2").Value = "VALIDACION"
$ws2.Range("F2").Value = "PUNTAJE"

# --- Rules table body (rows 3-11) ---
$ws2.Range("B3").Value = "DIUNSA"
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = "EDAD"
$ws2.Range("E3").Value = ">20"
$ws2.Range("F3").Value = 50

$ws2.Range("B4").Value = "DIUNSA"
$ws2.Range("C4").Value = 2
$ws2.Range("D4").Value = "EDAD"
$ws2.Range("E4").Value = ">25"
$ws2.Range("F4").Value = 15

$ws2.Range("B5").Value = "DIUNSA"
$ws2.Range("C5").Value = 3
$ws2.Range("D5").Value = "DEPENDIENTES"
$ws2.Range("E5").Value = 0
$ws2.Range("F5").Value = 20

$ws2.Range("B6").Value = "DIUNSA"
$ws2.Range("C6").Value = 4
$ws2.Range("D6").Value = "DEPENDIENTES"
$ws2.Range("E6").Value = ">2"
$ws2.Range("F6").Value = 10

$ws2.Range("B7").Value = "DIUNSA"
$ws2.Range("C7").Value = 5
$ws2.Range("D7").Value = "TRABAJO"
$ws2.Range("E7").Value = 15000
$ws2.Range("F7").Value = 15

$ws2.Range("B8").Value = "DIUNSA"
$ws2.Range("C8").Value = 6
$ws2.Range("D8").Value = "TRABAJO"
$ws2.Range("E8").Value = 20000
$ws2.Range("F8").Value = 25

$ws2.Range("B9").Value = "DIUNSA"
$ws2.Range("C9").Value = 7
$ws2.Range("D9").Value = "CREDITO"
$ws2.Range("E9").Value = "BUENO"
$ws2.Range("F9").Value = 20

$ws2.Range("B10").Value = "DIUNSA"
$ws2.Range("C10").Value = 8
$ws2.Range("D10").Value = "CREDITO"
$ws2.Range("E10").Value = "EXCELENTE"
$ws2.Range("F10").Value = 30

$ws2.Range("B11").Value = "DIUNSA"
$ws2.Range("C11").Value = 9
$ws2.Range("D11").Value = "OTRO"
$ws2.Range("E11").Value = "OTRO"
$ws2.Range("F11").Value = 30

# --- Second table: header (row 17) + sample scored record (row 18) ---
$ws2.Range("A17").Value = "fecah"
$ws2.Range("B17").Value = "SCORE"
$ws2.Range("C17").Value = "IDENTIDAD"
$ws2.Range("D17").Value = "EDAD"
$ws2.Range("E17").Value = "DEPENDIENTE"
$ws2.Range("F17").Value = "TRABAJO"
$ws2.Range("G17").Value = "CREDITO"
$ws2.Range("H17").Value = "OTRO"
$ws2.Range("I17").Value = "TOTAL SCORE"

$ws2.Range("A18").Value = ";l';;;;;;;;;;;;;;;"
$ws2.Range("B18").Value = 1
$ws2.Range("C18").Value = "ARMANDO"
$ws2.Range("D18").Value = 10
$ws2.Range("E18").Value = 20
$ws2.Range("F18").Value = 25
$ws2.Range("G18").Value = 30
$ws2.Range("H18").Value = 30
$ws2.Range("I18").Formula = "=SUM(D18:H18)"

# --- Column widths to match the published layout ---
$ws2.Columns.Item(4).ColumnWidth = 12.42578125
$ws2.Columns.Item(5).ColumnWidth = 15.7109375
$ws2.Columns.Item(8).ColumnWidth = 12

# --- Formatting ---
# Column E is centred throughout (plain style for most rows, currency for
# the two work-experience/income threshold rows).
$eRange = $ws2.Range("E2:E6,E9:E11,E17:E18")
$eRange.HorizontalAlignment = -4108

$eCurrency = $ws2.Range("E7:E8")
$eCurrency.HorizontalAlignment = -4108
$eCurrency.NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# A18 carries a date-formatted (but text) value in the source data.
$ws2.Range("A18").NumberFormat = "mm-dd-yy"

# Final selection on the new sheet, matching the saved workbook state.
$ws2.Range("A18").Select()
